$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = [ordered]@{
    "D2" = "29.271.23"
    "D3" = "1.898.50"
    "E3" = "  -0.48%  "
    "D4" = "1.001"
    "E4" = "  -0.35%  "
    "D5" = "325.82"
    "E5" = "  -0.53%  "
    "E6" = "  -0.41%  "
    "D7" = "0.4640"
    "E7" = "  -0.19%  "
    "D8" = "0.3921"
    "E8" = "  -0.35%  "
    "E9" = "  -0.87%  "
    "D10" = "0.9886"
    "E10" = "  -1.70%  "
    "E11" = "  -1.93%  "
    "D12" = "1.910.47"
    "E12" = "  +0.01%  "
    "D13" = "7.080"
    "E13" = "  -1.20%  "
    "D14" = "5.753"
    "E14" = "  -0.61%  "
    "D15" = "0.06994"
    "E15" = "  -0.05%  "
    "D16" = "88.46"
    "E16" = "  -0.32%  "
    "E17" = "  -0.41%  "
    "D18" = "0.00001001"
    "E18" = "  -0.96%  "
    "D19" = "17.10"
    "E19" = "  -0.97%  "
    "E20" = "  -0.32%  "
    "D21" = "29.284.71"
    "E21" = "  +0.41%  "
    "D22" = "5.294"
    "E22" = "  -1.76%  "
    "D23" = "11.10"
    "E23" = "  +0.18%  "
    "D24" = "2.095"
    "E24" = "  +2.25%  "
    "D25" = "156.13"
    "E25" = "  -0.39%  "
    "D26" = "19.43"
    "E26" = "  -0.85%  "
    "D27" = "6.028"
    "E27" = "  +2.41%  "
    "D28" = "118.59"
    "E28" = "  -0.95%  "
    "D29" = "1.924"
    "E29" = "  -4.35%  "
    "D30" = "0.09376"
    "E30" = "  -0.05%  "
    "D31" = "0.9054"
    "E31" = "  -2.26%  "
    "D32" = "5.299"
    "E32" = "  -1.32%  "
    "E33" = "  -1.40%  "
    "D34" = "3.229"
    "E34" = "  -1.57%  "
    "B35" = "TrustWalletToken"
    "C35" = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
    "D35" = "1.184"
    "E35" = "  +2.02%  "
    "B36" = "Hedera"
    "C36" = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
    "D36" = "0.05804"
    "E36" = "  -0.75%  "
    "D37" = "0.02092"
    "E37" = "  -0.64%  "
    "B38" = "FraxShare"
    "C38" = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
    "D38" = "7.803"
    "E38" = "  -2.66%  "
    "B39" = "Frax"
    "C39" = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
    "D39" = "0.9996"
    "E39" = "  -0.43%  "
    "D40" = "0.5714"
    "E40" = "  -0.99%  "
    "D41" = "0.1786"
    "E41" = "  -1.44%  "
    "D42" = "9.754"
    "E42" = "  -2.49%  "
    "E43" = "  -1.08%  "
    "D44" = "2.218"
    "E44" = "  -1.65%  "
    "D45" = "0.5359"
    "E45" = "  -1.51%  "
    "D46" = "0.07049"
    "E46" = "  -1.08%  "
    "D47" = "1.864"
    "E47" = "  -1.26%  "
    "D48" = "2.566"
    "E48" = "  +0.73%  "
    "D49" = "112.92"
    "E49" = "  +0.36%  "
    "D50" = "1.064"
    "E50" = "  -4.75%  "
    "D51" = "71.34"
    "E51" = "  -1.16%  "
}

foreach ($addr in $changes.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $changes[$addr]
    $cell.Style = "Normal"
}
